# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E23) listed the most recent periods in
# descending order (1704, 1703, 1702, 1701, 1612, 1611, 1610, 1609).
# The database was refreshed and the periods now run in ascending order
# starting from the new batch (1609, 1610, 1611, 1612, 1701, 1702, 1703,
# 1704), i.e. part 1 of the new account-statement periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1609"
$ws.Range("E17").Value = "1610"
$ws.Range("E18").Value = "1611"
$ws.Range("E19").Value = "1612"
$ws.Range("E20").Value = "1701"
$ws.Range("E21").Value = "1702"
$ws.Range("E22").Value = "1703"
$ws.Range("E23").Value = "1704"
